# Auto-generated Excel COM-interop edit script.
# Rebuilds 展览 (sheet1) rows 2-34 and 全部类型 (sheet4) rows 2-21 with the
# post-edit content, zeroes out 演出 (sheet2) F2:F6, and leaves 本地生活
# (sheet3) untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Make sure newly-extended rows inherit the same bordered/bold/centered
# style the existing index column (A) cells use.
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("A2:A34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$data1 = New-Object 'object[,]' 33,9
$data1[0,0] = 1
$data1[0,1] = '''2018-02-21'
$data1[0,2] = '2018南宁月邪动漫冬季盛典'
$data1[0,3] = '南宁国际会展中心  南宁国际会展中心'
$data1[0,4] = '2018.02.21 10:00-02.22 17:00'
$data1[0,5] = 0
$data1[0,6] = '已结束'
$data1[0,7] = 'https://show.bilibili.com/platform/detail.html?id=11460'
$data1[0,8] = '//i2.hdslb.com/bfs/openplatform/201712/immQq0qWGNybA.jpeg'
$data1[1,0] = 2
$data1[1,1] = '''2018-04-29'
$data1[1,2] = '2018年中国—东盟博览会动漫游戏展'
$data1[1,3] = '南宁国际会展中心  南宁国际会展中心'
$data1[1,4] = '2018.04.29 09:00-05.01 17:00'
$data1[1,5] = 0
$data1[1,6] = '已结束'
$data1[1,7] = 'https://show.bilibili.com/platform/detail.html?id=11931'
$data1[1,8] = '//i1.hdslb.com/bfs/openplatform/201801/imjtTxj0ZROoU.png'
$data1[2,0] = 3
$data1[2,1] = '''2018-07-28'
$data1[2,2] = '2018南宁月邪夏季盛典'
$data1[2,3] = '南宁国际会展中心  南宁国际会展中心'
$data1[2,4] = '2018.07.28 09:00-07.29 17:00'
$data1[2,5] = 0
$data1[2,6] = '已结束'
$data1[2,7] = 'https://show.bilibili.com/platform/detail.html?id=13412'
$data1[2,8] = '//i0.hdslb.com/bfs/openplatform/201806/imGkI6mLyvPgI.jpeg'
$data1[3,0] = 4
$data1[3,1] = '''2018-10-01'
$data1[3,2] = '南宁·2018Climax动漫游戏嘉年华'
$data1[3,3] = '南宁国际会展中心  南宁国际会展中心'
$data1[3,4] = '2018.10.01 10:00-10.02 17:00'
$data1[3,5] = 0
$data1[3,6] = '已结束'
$data1[3,7] = 'https://show.bilibili.com/platform/detail.html?id=14526'
$data1[3,8] = '//i0.hdslb.com/bfs/openplatform/201809/imRPAVGMfQk8Y.jpeg'
$data1[4,0] = 5
$data1[4,1] = '''2019-06-08'
$data1[4,2] = '南宁·微光国际动漫展'
$data1[4,3] = '沙井大道56号 华南城会展中心'
$data1[4,4] = '2019.06.08 10:00-06.09 17:00'
$data1[4,5] = 0
$data1[4,6] = '已结束'
$data1[4,7] = 'https://show.bilibili.com/platform/detail.html?id=17469'
$data1[4,8] = '//i2.hdslb.com/bfs/openplatform/201905/imXup7ya8VKis.jpeg'
$data1[5,0] = 6
$data1[5,1] = '''2021-10-30'
$data1[5,2] = '南宁·万圣漫控嘉年华04'
$data1[5,3] = '亭洪路45号17栋百益·上河城C7 百益上河城艺术中心'
$data1[5,4] = '2021.10.30 10:30-10.31 22:00'
$data1[5,5] = 0
$data1[5,6] = '已结束'
$data1[5,7] = 'https://show.bilibili.com/platform/detail.html?id=37560'
$data1[5,8] = '//i0.hdslb.com/bfs/openplatform/202110/FW9b4aQo1634800385093.jpeg'
$data1[6,0] = 7
$data1[6,1] = '''2022-01-01'
$data1[6,2] = '南宁·首届萌卡动漫嘉年华（取消）'
$data1[6,3] = '民族大道106号 南宁国际会展中心'
$data1[6,4] = '2022.01.01 09:00-01.02 18:00'
$data1[6,5] = 0
$data1[6,6] = '已结束'
$data1[6,7] = 'https://show.bilibili.com/platform/detail.html?id=54874'
$data1[6,8] = '//i2.hdslb.com/bfs/openplatform/202111/rMnkrbx11637654412833.jpeg'
$data1[7,0] = 8
$data1[7,1] = '''2022-03-26'
$data1[7,2] = '南宁·2022月邪动漫冬季盛典'
$data1[7,3] = '民族大道106号 南宁国际会展中心'
$data1[7,4] = '2022.03.26 09:30-03.27 17:30'
$data1[7,5] = 0
$data1[7,6] = '已结束'
$data1[7,7] = 'https://show.bilibili.com/platform/detail.html?id=55094'
$data1[7,8] = '//i0.hdslb.com/bfs/openplatform/202203/S1miyQb81647426988782.jpeg'
$data1[8,0] = 9
$data1[8,1] = '''2022-06-03'
$data1[8,2] = '【会员购严选】南宁·艾妮X漫控 潮流嘉年华05暨2022广西动漫文旅产业博览会'
$data1[8,3] = '亭洪路45号17栋百益·上河城C7 百益上河城艺术中心'
$data1[8,4] = '2022.06.03 09:30-06.05 17:00'
$data1[8,5] = 0
$data1[8,6] = '已结束'
$data1[8,7] = 'https://show.bilibili.com/platform/detail.html?id=58438'
$data1[8,8] = '//i1.hdslb.com/bfs/openplatform/202205/7005Fbvi1653537742378.jpeg'
$data1[9,0] = 10
$data1[9,1] = '''2023-07-22'
$data1[9,2] = '南宁·2023良牙动漫夏季盛典（夏典）'
$data1[9,3] = '民族大道106号 南宁国际会展中心'
$data1[9,4] = '2023.07.22 09:30-07.23 17:30'
$data1[9,5] = 0
$data1[9,6] = '已结束'
$data1[9,7] = 'https://show.bilibili.com/platform/detail.html?id=73723'
$data1[9,8] = '//i2.hdslb.com/bfs/openplatform/202306/dhGyvyqr1686648298409.jpeg'
$data1[10,0] = 11
$data1[10,1] = '''2023-10-28'
$data1[10,2] = '南宁·万圣漫控嘉年华08'
$data1[10,3] = '亭洪路45号 百益上河城'
$data1[10,4] = '2023.10.28 12:00-10.29 22:00'
$data1[10,5] = 0
$data1[10,6] = '已结束'
$data1[10,7] = 'https://show.bilibili.com/platform/detail.html?id=73274'
$data1[10,8] = '//i2.hdslb.com/bfs/openplatform/202310/7zvp4YhB1697698051810.jpeg'
$data1[11,0] = 12
$data1[11,1] = '''2023-12-09'
$data1[11,2] = '南宁·AP动漫游戏嘉年华内场票-倒霉死勒'
$data1[11,3] = '亭洪路45号 百益上河城'
$data1[11,4] = '2023.12.09 09:00-12.09 17:00'
$data1[11,5] = 0
$data1[11,6] = '已结束'
$data1[11,7] = 'https://show.bilibili.com/platform/detail.html?id=77715'
$data1[11,8] = '//i0.hdslb.com/bfs/openplatform/202310/NPLe3TrR1698288731028.jpeg'
$data1[12,0] = 13
$data1[12,1] = '''2024-02-15'
$data1[12,2] = '南宁·2024良牙动漫冬季盛典（冬典）'
$data1[12,3] = '民族大道106号 南宁国际会展中心'
$data1[12,4] = '2024.02.15 09:30-02.16 17:30'
$data1[12,5] = 0
$data1[12,6] = '已结束'
$data1[12,7] = 'https://show.bilibili.com/platform/detail.html?id=77938'
$data1[12,8] = '//i1.hdslb.com/bfs/openplatform/202311/YriBERx81701329557375.jpeg'
$data1[13,0] = 14
$data1[13,1] = '''2024-03-16'
$data1[13,2] = '南宁·草莓动漫节'
$data1[13,3] = '亭洪路45号 百益上河城'
$data1[13,4] = '2024.03.16 09:00-03.17 17:00'
$data1[13,5] = 0
$data1[13,6] = '已结束'
$data1[13,7] = 'https://show.bilibili.com/platform/detail.html?id=80943'
$data1[13,8] = '//i0.hdslb.com/bfs/openplatform/202402/vF9kexbx1707289709364.jpeg'
$data1[14,0] = 15
$data1[14,1] = '''2024-05-01'
$data1[14,2] = '南宁·2024三月三国潮动漫节（良牙春典）'
$data1[14,3] = '民族大道106号 南宁国际会展中心'
$data1[14,4] = '2024.05.01 09:30-05.02 17:30'
$data1[14,5] = 0
$data1[14,6] = '已结束'
$data1[14,7] = 'https://show.bilibili.com/platform/detail.html?id=82416'
$data1[14,8] = '//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg'
$data1[15,0] = 16
$data1[15,1] = '''2024-05-25'
$data1[15,2] = '南宁·第五人格Only1.0'
$data1[15,3] = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$data1[15,4] = '2024.05.25 10:00-05.25 17:30'
$data1[15,5] = 0
$data1[15,6] = '已结束'
$data1[15,7] = 'https://show.bilibili.com/platform/detail.html?id=84954'
$data1[15,8] = '//i0.hdslb.com/bfs/openplatform/202404/w5iZT4wE1714189905443.jpeg'
$data1[16,0] = 17
$data1[16,1] = '''2024-06-09'
$data1[16,2] = '南宁·AP动漫游戏嘉年华'
$data1[16,3] = '南宁国际会展中心  南宁国际会展中心'
$data1[16,4] = '2024.06.09 09:00-06.10 17:00'
$data1[16,5] = 0
$data1[16,6] = '已结束'
$data1[16,7] = 'https://show.bilibili.com/platform/detail.html?id=84793'
$data1[16,8] = '//i1.hdslb.com/bfs/openplatform/202405/hyC2ZhnZ1715826721453.jpeg'
$data1[17,0] = 18
$data1[17,1] = '''2024-06-09'
$data1[17,2] = '宾阳·荷止国风动漫展'
$data1[17,3] = '商贸城社区南段86-1 宾阳金玉酒店'
$data1[17,4] = '2024.06.09 10:30-06.09 16:30'
$data1[17,5] = 0
$data1[17,6] = '已结束'
$data1[17,7] = 'https://show.bilibili.com/platform/detail.html?id=85980'
$data1[17,8] = '//i2.hdslb.com/bfs/openplatform/202405/EhUqTg5l1715838043315.jpeg'
$data1[18,0] = 19
$data1[18,1] = '''2024-06-15'
$data1[18,2] = '南宁·星STAR国潮嘉年华（取消）'
$data1[18,3] = '亭洪路45号 百益上河城'
$data1[18,4] = '2024.06.15 09:00-06.16 17:00'
$data1[18,5] = 0
$data1[18,6] = '已结束'
$data1[18,7] = 'https://show.bilibili.com/platform/detail.html?id=86198'
$data1[18,8] = '//i0.hdslb.com/bfs/openplatform/202405/orwMgait1716448294056.jpeg'
$data1[19,0] = 20
$data1[19,1] = '''2024-06-22'
$data1[19,2] = '南宁·排球少年ONLY（取消）'
$data1[19,3] = '亭洪路45号 水明漾宴会中心'
$data1[19,4] = '2024.06.22 09:45-06.22 17:00'
$data1[19,5] = 0
$data1[19,6] = '已结束'
$data1[19,7] = 'https://show.bilibili.com/platform/detail.html?id=86465'
$data1[19,8] = '//i0.hdslb.com/bfs/openplatform/202405/GaaD97dL1716883956953.jpeg'
$data1[20,0] = 21
$data1[20,1] = '''2024-07-06'
$data1[20,2] = '南宁·小蜜蜂动漫嘉年华2.0'
$data1[20,3] = '亭洪路45号 百益上河城'
$data1[20,4] = '2024.07.06 10:00-07.06 17:00'
$data1[20,5] = 0
$data1[20,6] = 55
$data1[20,7] = 'https://show.bilibili.com/platform/detail.html?id=84925'
$data1[20,8] = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$data1[21,0] = 22
$data1[21,1] = '''2024-07-06'
$data1[21,2] = '南宁·首届童话梦境Lolita茶会'
$data1[21,3] = '明秀东路157号 利泰国际大酒店'
$data1[21,4] = '2024.07.06 13:00-07.06 17:00'
$data1[21,5] = 0
$data1[21,6] = '已停售'
$data1[21,7] = 'https://show.bilibili.com/platform/detail.html?id=85776'
$data1[21,8] = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$data1[22,0] = 23
$data1[22,1] = '''2024-07-12'
$data1[22,2] = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$data1[22,3] = '民族大道106号 南宁国际会展中心'
$data1[22,4] = '2024.07.12 09:30-07.14 17:00'
$data1[22,5] = 0
$data1[22,6] = 50
$data1[22,7] = 'https://show.bilibili.com/platform/detail.html?id=87182'
$data1[22,8] = '//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg'
$data1[23,0] = 24
$data1[23,1] = '''2024-07-13'
$data1[23,2] = '南宁·0713国乙ONLY'
$data1[23,3] = '亭洪路45号 水明漾宴会中心'
$data1[23,4] = '2024.07.13 09:30-07.13 21:00'
$data1[23,5] = 0
$data1[23,6] = 68
$data1[23,7] = 'https://show.bilibili.com/platform/detail.html?id=86378'
$data1[23,8] = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'
$data1[24,0] = 25
$data1[24,1] = '''2024-07-14'
$data1[24,2] = '广西·首届明日方舟only展 - 花庭圣梦'
$data1[24,3] = '明秀东路157号 利泰国际大酒店'
$data1[24,4] = '2024.07.14 09:00-07.14 18:00'
$data1[24,5] = 0
$data1[24,6] = 69
$data1[24,7] = 'https://show.bilibili.com/platform/detail.html?id=85852'
$data1[24,8] = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'
$data1[25,0] = 26
$data1[25,1] = '''2024-07-20'
$data1[25,2] = '南宁·AB动漫游戏嘉年华'
$data1[25,3] = '五象大道西段669号 广西体育中心体育馆'
$data1[25,4] = '2024.07.20 09:30-07.21 17:00'
$data1[25,5] = 0
$data1[25,6] = 60
$data1[25,7] = 'https://show.bilibili.com/platform/detail.html?id=84862'
$data1[25,8] = '//i1.hdslb.com/bfs/openplatform/202407/R7iP9Iio1720170437964.jpeg'
$data1[26,0] = 27
$data1[26,1] = '''2024-07-20'
$data1[26,2] = '横州·第二届海棠动漫游戏嘉年华'
$data1[26,3] = '茉莉花大道 横州国际大酒店'
$data1[26,4] = '2024.07.20 09:30-07.20 17:00'
$data1[26,5] = 0
$data1[26,6] = 30
$data1[26,7] = 'https://show.bilibili.com/platform/detail.html?id=84799'
$data1[26,8] = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'
$data1[27,0] = 28
$data1[27,1] = '''2024-07-27'
$data1[27,2] = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$data1[27,3] = '民族大道106号 南宁国际会展中心'
$data1[27,4] = '2024.07.27 09:30-07.28 17:30'
$data1[27,5] = 0
$data1[27,6] = 55
$data1[27,7] = 'https://show.bilibili.com/platform/detail.html?id=85264'
$data1[27,8] = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'
$data1[28,0] = 29
$data1[28,1] = '''2024-08-03'
$data1[28,2] = '南宁·火影忍者only'
$data1[28,3] = '厢竹大道65号 桔子酒店'
$data1[28,4] = '2024.08.03 10:00-08.03 17:00'
$data1[28,5] = 0
$data1[28,6] = 68
$data1[28,7] = 'https://show.bilibili.com/platform/detail.html?id=86994'
$data1[28,8] = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$data1[29,0] = 30
$data1[29,1] = '''2024-08-03'
$data1[29,2] = '南宁·蔚蓝档案only'
$data1[29,3] = '亭洪路45号 百益上河城'
$data1[29,4] = '2024.08.03 09:00-08.03 17:00'
$data1[29,5] = 0
$data1[29,6] = 68
$data1[29,7] = 'https://show.bilibili.com/platform/detail.html?id=85370'
$data1[29,8] = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'
$data1[30,0] = 31
$data1[30,1] = '''2024-08-10'
$data1[30,2] = '南宁·国乙only'
$data1[30,3] = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$data1[30,4] = '2024.08.10 10:00-08.10 17:00'
$data1[30,5] = 0
$data1[30,6] = 40
$data1[30,7] = 'https://show.bilibili.com/platform/detail.html?id=88227'
$data1[30,8] = '//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg'
$data1[31,0] = 32
$data1[31,1] = '''2024-08-24'
$data1[31,2] = '南宁·第二届北极光动漫展'
$data1[31,3] = '民族大道106号 南宁国际会展中心'
$data1[31,4] = '2024.08.24 09:00-08.25 17:00'
$data1[31,5] = 0
$data1[31,6] = 65
$data1[31,7] = 'https://show.bilibili.com/platform/detail.html?id=88276'
$data1[31,8] = '//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg'
$data1[32,0] = 33
$data1[32,1] = '''2024-11-02'
$data1[32,2] = '南宁·万圣漫控嘉年华10'
$data1[32,3] = '亭洪路45号 百益上河城'
$data1[32,4] = '2024.11.02 11:00-11.03 22:00'
$data1[32,5] = 0
$data1[32,6] = 50
$data1[32,7] = 'https://show.bilibili.com/platform/detail.html?id=87820'
$data1[32,8] = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'
$ws1.Range("A2:I34").Value = $data1

# ---------------------------------------------------------------
# Sheet "演出" (performances) -- only the "want to go" counter resets
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("F6").Value = 0

# ---------------------------------------------------------------
# Sheet "本地生活" (local life) -- unchanged, nothing to do.
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("A1").Copy() | Out-Null
$ws4.Range("A2:A21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$data4 = New-Object 'object[,]' 20,9
$data4[0,0] = 1
$data4[0,1] = '''2024-07-06'
$data4[0,2] = '南宁·小蜜蜂动漫嘉年华2.0'
$data4[0,3] = '亭洪路45号 百益上河城'
$data4[0,4] = '2024.07.06 10:00-07.06 17:00'
$data4[0,5] = 0
$data4[0,6] = 55
$data4[0,7] = 'https://show.bilibili.com/platform/detail.html?id=84925'
$data4[0,8] = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$data4[1,0] = 2
$data4[1,1] = '''2024-07-06'
$data4[1,2] = '南宁·首届童话梦境Lolita茶会'
$data4[1,3] = '明秀东路157号 利泰国际大酒店'
$data4[1,4] = '2024.07.06 13:00-07.06 17:00'
$data4[1,5] = 0
$data4[1,6] = '已停售'
$data4[1,7] = 'https://show.bilibili.com/platform/detail.html?id=85776'
$data4[1,8] = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$data4[2,0] = 3
$data4[2,1] = '''2024-07-06'
$data4[2,2] = '南宁·首届童话梦境Lolita茶会'
$data4[2,3] = '明秀东路157号 利泰国际大酒店'
$data4[2,4] = '2024.07.06 13:00-07.06 17:00'
$data4[2,5] = 0
$data4[2,6] = '已停售'
$data4[2,7] = 'https://show.bilibili.com/platform/detail.html?id=85776'
$data4[2,8] = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$data4[3,0] = 4
$data4[3,1] = '''2024-07-12'
$data4[3,2] = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$data4[3,3] = '民族大道106号 南宁国际会展中心'
$data4[3,4] = '2024.07.12 09:30-07.14 17:00'
$data4[3,5] = 0
$data4[3,6] = 50
$data4[3,7] = 'https://show.bilibili.com/platform/detail.html?id=87182'
$data4[3,8] = '//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg'
$data4[4,0] = 5
$data4[4,1] = '''2024-07-13'
$data4[4,2] = '南宁·0713国乙ONLY'
$data4[4,3] = '亭洪路45号 水明漾宴会中心'
$data4[4,4] = '2024.07.13 09:30-07.13 21:00'
$data4[4,5] = 0
$data4[4,6] = 68
$data4[4,7] = 'https://show.bilibili.com/platform/detail.html?id=86378'
$data4[4,8] = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'
$data4[5,0] = 6
$data4[5,1] = '''2024-07-14'
$data4[5,2] = '广西·首届明日方舟only展 - 花庭圣梦'
$data4[5,3] = '明秀东路157号 利泰国际大酒店'
$data4[5,4] = '2024.07.14 09:00-07.14 18:00'
$data4[5,5] = 0
$data4[5,6] = 69
$data4[5,7] = 'https://show.bilibili.com/platform/detail.html?id=85852'
$data4[5,8] = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'
$data4[6,0] = 7
$data4[6,1] = '''2024-07-18'
$data4[6,2] = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》'
$data4[6,3] = '龙堤路25号 广西文化艺术中心'
$data4[6,4] = '2024.07.18 20:00-07.18 21:30'
$data4[6,5] = 0
$data4[6,6] = 108
$data4[6,7] = 'https://show.bilibili.com/platform/detail.html?id=85816'
$data4[6,8] = '//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg'
$data4[7,0] = 8
$data4[7,1] = '''2024-07-19'
$data4[7,2] = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 '
$data4[7,3] = '龙堤路25号 广西文化艺术中心'
$data4[7,4] = '2024.07.19 20:00-07.19 22:00'
$data4[7,5] = 0
$data4[7,6] = 108
$data4[7,7] = 'https://show.bilibili.com/platform/detail.html?id=85831'
$data4[7,8] = '//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg'
$data4[8,0] = 9
$data4[8,1] = '''2024-07-20'
$data4[8,2] = '南宁·AB动漫游戏嘉年华'
$data4[8,3] = '五象大道西段669号 广西体育中心体育馆'
$data4[8,4] = '2024.07.20 09:30-07.21 17:00'
$data4[8,5] = 0
$data4[8,6] = 60
$data4[8,7] = 'https://show.bilibili.com/platform/detail.html?id=84862'
$data4[8,8] = '//i1.hdslb.com/bfs/openplatform/202407/R7iP9Iio1720170437964.jpeg'
$data4[9,0] = 10
$data4[9,1] = '''2024-07-20'
$data4[9,2] = '横州·第二届海棠动漫游戏嘉年华'
$data4[9,3] = '茉莉花大道 横州国际大酒店'
$data4[9,4] = '2024.07.20 09:30-07.20 17:00'
$data4[9,5] = 0
$data4[9,6] = 30
$data4[9,7] = 'https://show.bilibili.com/platform/detail.html?id=84799'
$data4[9,8] = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'
$data4[10,0] = 11
$data4[10,1] = '''2024-07-21'
$data4[10,2] = '南宁·跨越二次元ACG神级动漫世界巡回演唱会——'
$data4[10,3] = '中山路万象汇L2层37号 候朋现场HOPELIVE-中山路万象汇店'
$data4[10,4] = '2024.07.21 20:00-07.21 22:00'
$data4[10,5] = 0
$data4[10,6] = 138
$data4[10,7] = 'https://show.bilibili.com/platform/detail.html?id=88699'
$data4[10,8] = '//i1.hdslb.com/bfs/openplatform/202407/uvwreDk61720071220041.jpeg'
$data4[11,0] = 12
$data4[11,1] = '''2024-07-27'
$data4[11,2] = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$data4[11,3] = '民族大道106号 南宁国际会展中心'
$data4[11,4] = '2024.07.27 09:30-07.28 17:30'
$data4[11,5] = 0
$data4[11,6] = 55
$data4[11,7] = 'https://show.bilibili.com/platform/detail.html?id=85264'
$data4[11,8] = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'
$data4[12,0] = 13
$data4[12,1] = '''2024-08-03'
$data4[12,2] = '南宁·火影忍者only'
$data4[12,3] = '厢竹大道65号 桔子酒店'
$data4[12,4] = '2024.08.03 10:00-08.03 17:00'
$data4[12,5] = 0
$data4[12,6] = 68
$data4[12,7] = 'https://show.bilibili.com/platform/detail.html?id=86994'
$data4[12,8] = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$data4[13,0] = 14
$data4[13,1] = '''2024-08-03'
$data4[13,2] = '南宁·火影忍者only'
$data4[13,3] = '厢竹大道65号 桔子酒店'
$data4[13,4] = '2024.08.03 10:00-08.03 17:00'
$data4[13,5] = 0
$data4[13,6] = 68
$data4[13,7] = 'https://show.bilibili.com/platform/detail.html?id=86994'
$data4[13,8] = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$data4[14,0] = 15
$data4[14,1] = '''2024-08-03'
$data4[14,2] = '南宁·蔚蓝档案only'
$data4[14,3] = '亭洪路45号 百益上河城'
$data4[14,4] = '2024.08.03 09:00-08.03 17:00'
$data4[14,5] = 0
$data4[14,6] = 68
$data4[14,7] = 'https://show.bilibili.com/platform/detail.html?id=85370'
$data4[14,8] = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'
$data4[15,0] = 16
$data4[15,1] = '''2024-08-10'
$data4[15,2] = '南宁·国乙only'
$data4[15,3] = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$data4[15,4] = '2024.08.10 10:00-08.10 17:00'
$data4[15,5] = 0
$data4[15,6] = 40
$data4[15,7] = 'https://show.bilibili.com/platform/detail.html?id=88227'
$data4[15,8] = '//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg'
$data4[16,0] = 17
$data4[16,1] = '''2024-08-10'
$data4[16,2] = '南宁·限时7折|浪漫七夕《一生所爱》《爱乐之城》《假如爱有天意》经典浪漫电影主题音乐会'
$data4[16,3] = '龙堤路25号 广西文化艺术中心'
$data4[16,4] = '2024.08.10 20:00-08.10 21:30'
$data4[16,5] = 0
$data4[16,6] = 99
$data4[16,7] = 'https://show.bilibili.com/platform/detail.html?id=87729'
$data4[16,8] = '//i1.hdslb.com/bfs/openplatform/202406/qKUDMYOh1718177639735.png'
$data4[17,0] = 18
$data4[17,1] = '''2024-08-14'
$data4[17,2] = '南宁·新西兰治愈系民谣歌手LukeThompson2024中国巡演 KEEP ROLLING ON '
$data4[17,3] = '中山路万象汇L2层37号 候朋现场HOPELIVE-中山路万象汇店'
$data4[17,4] = '2024.08.14 20:00-08.14 21:30'
$data4[17,5] = 0
$data4[17,6] = 180
$data4[17,7] = 'https://show.bilibili.com/platform/detail.html?id=88015'
$data4[17,8] = '//i1.hdslb.com/bfs/openplatform/202406/76WI4tA01718179482365.jpeg'
$data4[18,0] = 19
$data4[18,1] = '''2024-08-24'
$data4[18,2] = '南宁·第二届北极光动漫展'
$data4[18,3] = '民族大道106号 南宁国际会展中心'
$data4[18,4] = '2024.08.24 09:00-08.25 17:00'
$data4[18,5] = 0
$data4[18,6] = 65
$data4[18,7] = 'https://show.bilibili.com/platform/detail.html?id=88276'
$data4[18,8] = '//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg'
$data4[19,0] = 20
$data4[19,1] = '''2024-11-02'
$data4[19,2] = '南宁·万圣漫控嘉年华10'
$data4[19,3] = '亭洪路45号 百益上河城'
$data4[19,4] = '2024.11.02 11:00-11.03 22:00'
$data4[19,5] = 0
$data4[19,6] = 50
$data4[19,7] = 'https://show.bilibili.com/platform/detail.html?id=87820'
$data4[19,8] = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'
$ws4.Range("A2:I21").Value = $data4
